# Saudi Arabia Division 1 - base update (25-05-2024 15:10)
# The update swaps the full record (all columns except the running index in column A)
# between specific pairs of adjacent rows in the sheet. Each pair represents two fixtures
# played on the same date that need to exchange places in the listing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (row1, row2) pairs whose B:AB content must be swapped
$pairs = @(
    @(2, 3),
    @(29, 30),
    @(35, 36),
    @(38, 39),
    @(68, 69),
    @(107, 108),
    @(172, 173),
    @(225, 226),
    @(230, 231),
    @(243, 244),
    @(259, 260),
    @(262, 263),
    @(280, 281),
    @(290, 291),
    @(292, 293)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AB$r1")
    $range2 = $ws.Range("B$r2`:AB$r2")

    $val1 = $range1.Value2
    $val2 = $range2.Value2

    $range1.Value2 = $val2
    $range2.Value2 = $val1
}
